# Add the "no-ras" run's throughput/ops-counts to the overview sheet.
# Fills in the previously-empty AL/AM columns (raw count + relative stddev)
# for each benchmark row; the dependent AN column (AL/$F ratio) recalculates
# automatically from its existing formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AL7").Value  = 2560
$ws.Range("AM7").Value  = 0.693

$ws.Range("AL8").Value  = 2383
$ws.Range("AM8").Value  = 1.67

$ws.Range("AL9").Value  = 3061
$ws.Range("AM9").Value  = 1.54

$ws.Range("AL10").Value = 2396
$ws.Range("AM10").Value = 0.681

$ws.Range("AL11").Value = 1399
$ws.Range("AM11").Value = 1.01

$ws.Range("AL12").Value = 2906
$ws.Range("AM12").Value = 0.607

$ws.Range("AL13").Value = 2337
$ws.Range("AM13").Value = 0.613

$ws.Range("AL14").Value = 2796
$ws.Range("AM14").Value = 0.61

$ws.Range("AL15").Value = 1455
$ws.Range("AM15").Value = 2.02

$ws.Range("AL16").Value = 7620
$ws.Range("AM16").Value = 0.811

# Row 18 only gets the relative-stddev summary value; AM18 stays blank.
$ws.Range("AL18").Value = 0.9251

# Update the view: zoomed in further, scrolled back to the top, and the
# active selection moved from AL24 to AH2.
$ws.Select()
$excel.ActiveWindow.Zoom = 150
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("AH2").Select()
